$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new header "VIMMP_DEF" in F1, matching the formatting/style of E1
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("F1").Value = "VIMMP_DEF"
